# Auto-generated script to update cached market price/profit values
# in the Leve profit tables across all 8 job sheets, per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value2 = 1616.65
$ws.Range("I2").Value2 = 322.3
$ws.Range("J2").Value2 = 2911
$ws.Range("K2").Value2 = 322.3
$ws.Range("L2").Value2 = 2911
$ws.Range("M2").Value2 = -209.3
$ws.Range("N2").Value2 = -3137
$ws.Range("H17").Value2 = 3202.4614
$ws.Range("J17").Value2 = 3202.4614
$ws.Range("L17").Value2 = 9607.3842
$ws.Range("N17").Value2 = -9943.3842
$ws.Range("H28").Value2 = 160.45454
$ws.Range("I28").Value2 = 159.3
$ws.Range("J28").Value2 = 172
$ws.Range("K28").Value2 = 159.3
$ws.Range("L28").Value2 = 172
$ws.Range("M28").Value2 = 325.7
$ws.Range("N28").Value2 = -1142
$ws.Range("H40").Value2 = 11557.143
$ws.Range("J40").Value2 = 11333.333
$ws.Range("L40").Value2 = 11333.333
$ws.Range("N40").Value2 = -11683.333
$ws.Range("H92").Value2 = 603.5833
$ws.Range("I92").Value2 = 618.8182
$ws.Range("K92").Value2 = 618.8182
$ws.Range("M92").Value2 = 629.1818
$ws.Range("H98").Value2 = 1322.1538
$ws.Range("I98").Value2 = 1393.1111
$ws.Range("K98").Value2 = 1393.1111
$ws.Range("M98").Value2 = 104.8888999999999
$ws.Range("H112").Value2 = 3250.8462
$ws.Range("J112").Value2 = 3328.88
$ws.Range("L112").Value2 = 9986.639999999999
$ws.Range("N112").Value2 = -12202.64
$ws.Range("H122").Value2 = 1322.1538
$ws.Range("I122").Value2 = 1393.1111
$ws.Range("K122").Value2 = 4179.3333
$ws.Range("M122").Value2 = -1729.3333
$ws.Range("H125").Value2 = 821.2
$ws.Range("J125").Value2 = 1107
$ws.Range("L125").Value2 = 9963
$ws.Range("N125").Value2 = -14883
$ws.Range("H127").Value2 = 999.5
$ws.Range("I127").Value2 = 999
$ws.Range("K127").Value2 = 2997
$ws.Range("M127").Value2 = 1963
$ws.Range("H129").Value2 = 1466.9333
$ws.Range("I129").Value2 = 467.66666
$ws.Range("J129").Value2 = 2965.8333
$ws.Range("K129").Value2 = 1402.99998
$ws.Range("L129").Value2 = 8897.499899999999
$ws.Range("M129").Value2 = 3597.00002
$ws.Range("N129").Value2 = -18897.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value2 = 626
$ws.Range("I5").Value2 = 490.22223
$ws.Range("J5").Value2 = 1033.3334
$ws.Range("K5").Value2 = 490.22223
$ws.Range("L5").Value2 = 1033.3334
$ws.Range("M5").Value2 = -378.22223
$ws.Range("N5").Value2 = -1257.3334
$ws.Range("H61").Value2 = 3001.476
$ws.Range("I61").Value2 = 2339.6667
$ws.Range("K61").Value2 = 2339.6667
$ws.Range("M61").Value2 = -2127.6667
$ws.Range("H88").Value2 = 1765.2354
$ws.Range("I88").Value2 = 1417.1111
$ws.Range("J88").Value2 = 2156.875
$ws.Range("K88").Value2 = 1417.1111
$ws.Range("L88").Value2 = 2156.875
$ws.Range("M88").Value2 = -1011.1111
$ws.Range("N88").Value2 = -2968.875
$ws.Range("H91").Value2 = 1765.2354
$ws.Range("I91").Value2 = 1417.1111
$ws.Range("J91").Value2 = 2156.875
$ws.Range("K91").Value2 = 1417.1111
$ws.Range("L91").Value2 = 2156.875
$ws.Range("M91").Value2 = -13.11110000000008
$ws.Range("N91").Value2 = -4964.875
$ws.Range("H136").Value2 = 3001.476
$ws.Range("I136").Value2 = 2339.6667
$ws.Range("K136").Value2 = 7019.000100000001
$ws.Range("M136").Value2 = -4469.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value2 = 626
$ws.Range("I4").Value2 = 490.22223
$ws.Range("J4").Value2 = 1033.3334
$ws.Range("K4").Value2 = 490.22223
$ws.Range("L4").Value2 = 1033.3334
$ws.Range("M4").Value2 = -375.22223
$ws.Range("N4").Value2 = -1263.3334
$ws.Range("H105").Value2 = 3857.25
$ws.Range("I105").Value2 = 3723
$ws.Range("K105").Value2 = 3723
$ws.Range("M105").Value2 = -1976

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value2 = 66788.39999999999
$ws.Range("I7").Value2 = 166796.67
$ws.Range("J7").Value2 = 116.22222
$ws.Range("K7").Value2 = 166796.67
$ws.Range("L7").Value2 = 116.22222
$ws.Range("M7").Value2 = -166683.67
$ws.Range("N7").Value2 = -342.22222
$ws.Range("H122").Value2 = 2304.818
$ws.Range("I122").Value2 = 2610.7646
$ws.Range("K122").Value2 = 7832.293799999999
$ws.Range("M122").Value2 = -5382.293799999999
$ws.Range("H141").Value2 = 171473.61
$ws.Range("J141").Value2 = 171473.61
$ws.Range("L141").Value2 = 171473.61
$ws.Range("N141").Value2 = -181833.61

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value2 = 190689.36
$ws.Range("I104").Value2 = 1819.2727
$ws.Range("J104").Value2 = 339087.28
$ws.Range("K104").Value2 = 5457.8181
$ws.Range("L104").Value2 = 1017261.84
$ws.Range("M104").Value2 = -2836.8181
$ws.Range("N104").Value2 = -1022503.84

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 7126.636
$ws.Range("I70").Value2 = 5224.5
$ws.Range("J70").Value2 = 8213.571
$ws.Range("K70").Value2 = 5224.5
$ws.Range("L70").Value2 = 8213.571
$ws.Range("M70").Value2 = -4954.5
$ws.Range("N70").Value2 = -8753.571
$ws.Range("H73").Value2 = 7126.636
$ws.Range("I73").Value2 = 5224.5
$ws.Range("J73").Value2 = 8213.571
$ws.Range("K73").Value2 = 5224.5
$ws.Range("L73").Value2 = 8213.571
$ws.Range("M73").Value2 = -4288.5
$ws.Range("N73").Value2 = -10085.571
$ws.Range("H80").Value2 = 10395.2
$ws.Range("I80").Value2 = 5494
$ws.Range("J80").Value2 = 30000
$ws.Range("K80").Value2 = 5494
$ws.Range("L80").Value2 = 30000
$ws.Range("M80").Value2 = -4496
$ws.Range("N80").Value2 = -31996
$ws.Range("H83").Value2 = 10395.2
$ws.Range("I83").Value2 = 5494
$ws.Range("J83").Value2 = 30000
$ws.Range("K83").Value2 = 27470
$ws.Range("L83").Value2 = 150000
$ws.Range("M83").Value2 = -22478
$ws.Range("N83").Value2 = -159984
$ws.Range("H102").Value2 = 3441.7917
$ws.Range("I102").Value2 = 2901.375
$ws.Range("K102").Value2 = 2901.375
$ws.Range("M102").Value2 = -1279.375
$ws.Range("H123").Value2 = 45735.176
$ws.Range("J123").Value2 = 43373.668
$ws.Range("L123").Value2 = 43373.668
$ws.Range("N123").Value2 = -48273.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value2 = 403.83334
$ws.Range("I9").Value2 = 384.6
$ws.Range("J9").Value2 = 500
$ws.Range("K9").Value2 = 384.6
$ws.Range("L9").Value2 = 500
$ws.Range("M9").Value2 = -160.6
$ws.Range("N9").Value2 = -948
$ws.Range("H22").Value2 = 739.2
$ws.Range("I22").Value2 = 739.2
$ws.Range("J22").Value2 = 0
$ws.Range("K22").Value2 = 739.2
$ws.Range("L22").Value2 = 0
$ws.Range("M22").Value2 = -444.2
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value2 = 739.2
$ws.Range("I27").Value2 = 739.2
$ws.Range("J27").Value2 = 0
$ws.Range("K27").Value2 = 739.2
$ws.Range("L27").Value2 = 0
$ws.Range("M27").Value2 = -632.2
$ws.Range("N27").ClearContents()
$ws.Range("H46").Value2 = 1878.3636
$ws.Range("I46").Value2 = 1558.8
$ws.Range("J46").Value2 = 2144.6667
$ws.Range("K46").Value2 = 1558.8
$ws.Range("L46").Value2 = 2144.6667
$ws.Range("M46").Value2 = -1370.8
$ws.Range("N46").Value2 = -2520.6667
$ws.Range("H82").Value2 = 2077.3333
$ws.Range("I82").Value2 = 2077.3333
$ws.Range("K82").Value2 = 2077.3333
$ws.Range("M82").Value2 = -1716.3333
$ws.Range("H85").Value2 = 2077.3333
$ws.Range("I85").Value2 = 2077.3333
$ws.Range("K85").Value2 = 2077.3333
$ws.Range("M85").Value2 = -829.3332999999998
$ws.Range("H132").Value2 = 3226.125
$ws.Range("I132").Value2 = 3398.25
$ws.Range("J132").Value2 = 3054
$ws.Range("K132").Value2 = 10194.75
$ws.Range("L132").Value2 = 9162
$ws.Range("M132").Value2 = -7664.75
$ws.Range("N132").Value2 = -14222
$ws.Range("H136").Value2 = 4177.35
$ws.Range("J136").Value2 = 0
$ws.Range("L136").Value2 = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value2 = 17929.9
$ws.Range("I62").Value2 = 15500
$ws.Range("J62").Value2 = 18537.375
$ws.Range("K62").Value2 = 15500
$ws.Range("L62").Value2 = 18537.375
$ws.Range("M62").Value2 = -14876
$ws.Range("N62").Value2 = -19785.375
$ws.Range("H65").Value2 = 17929.9
$ws.Range("I65").Value2 = 15500
$ws.Range("J65").Value2 = 18537.375
$ws.Range("K65").Value2 = 77500
$ws.Range("L65").Value2 = 92686.875
$ws.Range("M65").Value2 = -74380
$ws.Range("N65").Value2 = -98926.875
$ws.Range("H113").Value2 = 879.8
$ws.Range("I113").Value2 = 879.8
$ws.Range("K113").Value2 = 2639.4
$ws.Range("M113").Value2 = -469.3999999999996
$ws.Range("H126").Value2 = 2438.5
$ws.Range("I126").Value2 = 2259.9443
$ws.Range("J126").Value2 = 4045.5
$ws.Range("K126").Value2 = 6779.8329
$ws.Range("L126").Value2 = 12136.5
$ws.Range("M126").Value2 = -4309.8329
$ws.Range("N126").Value2 = -17076.5

